$d = $word.ActiveDocument

$target = "Added hover and active colours to the anchors in the nav bar."
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`n")
    if ($t -eq $target) {
        $r = $p.Range
        $r.InsertParagraphAfter()
        break
    }
}

$d.Paragraphs.Last.Range.Text = "Noticed text was fine on laptop and responsive view, but very small on mobile. So attempted to utilise further for mobile."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Made nav bar sticky, so it gets moved up until just under the header bar and becomes fixed."

Write-Output "done"
